$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")
$totaux = $wb.Worksheets.Item("Totaux")

# xlPasteFormats constant
$xlPasteFormats = -4122
# xlCenter constant (not needed directly, using format copy instead)

# --- Add the 5 new journal rows (48-52), copying number/alignment formats
#     from the row directly above so the cell styles match the existing
#     table rows exactly. ---

# Row 48: Entretien avec le responsable infrastructure
$ws.Range("A48").Value = 44995
$ws.Range("A47").Copy()
$ws.Range("A48").PasteSpecial($xlPasteFormats)
$ws.Range("B48").Value = 5
$ws.Range("B47").Copy()
$ws.Range("B48").PasteSpecial($xlPasteFormats)
$ws.Range("C48").Value = 0.03125
$ws.Range("C47").Copy()
$ws.Range("C48").PasteSpecial($xlPasteFormats)
$ws.Range("D48").Value = "Entretien"
$ws.Range("D47").Copy()
$ws.Range("D48").PasteSpecial($xlPasteFormats)
$ws.Range("E48").Value = "Entretien avec le responsable infrastructure"

# Row 49: Rédaction d'un rapport d'entretien
$ws.Range("A49").Value = 44995
$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial($xlPasteFormats)
$ws.Range("B49").Value = 5
$ws.Range("B48").Copy()
$ws.Range("B49").PasteSpecial($xlPasteFormats)
$ws.Range("C49").Value = 0.020833333333333332
$ws.Range("C48").Copy()
$ws.Range("C49").PasteSpecial($xlPasteFormats)
$ws.Range("D49").Value = "Documentation"
$ws.Range("D48").Copy()
$ws.Range("D49").PasteSpecial($xlPasteFormats)
$ws.Range("E49").Value = "Rédaction d'un rapport d'entretien "

# Row 50: Entretien avec l'équipe de développement
$ws.Range("A50").Value = 44995
$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial($xlPasteFormats)
$ws.Range("B50").Value = 5
$ws.Range("B49").Copy()
$ws.Range("B50").PasteSpecial($xlPasteFormats)
$ws.Range("C50").Value = 0.020833333333333332
$ws.Range("C49").Copy()
$ws.Range("C50").PasteSpecial($xlPasteFormats)
$ws.Range("D50").Value = "Entretien"
$ws.Range("D49").Copy()
$ws.Range("D50").PasteSpecial($xlPasteFormats)
$ws.Range("E50").Value = "Entretien avec l'équipe de développement "

# Row 51: Rédaction d'un rapport d'entretien (2nd occurrence)
$ws.Range("A51").Value = 44995
$ws.Range("A50").Copy()
$ws.Range("A51").PasteSpecial($xlPasteFormats)
$ws.Range("B51").Value = 5
$ws.Range("B50").Copy()
$ws.Range("B51").PasteSpecial($xlPasteFormats)
$ws.Range("C51").Value = 0.020833333333333332
$ws.Range("C50").Copy()
$ws.Range("C51").PasteSpecial($xlPasteFormats)
$ws.Range("D51").Value = "Documentation"
$ws.Range("D50").Copy()
$ws.Range("D51").PasteSpecial($xlPasteFormats)
$ws.Range("E51").Value = "Rédaction d'un rapport d'entretien "

# Row 52: Configuration des catégories pour le service Développement
$ws.Range("A52").Value = 44995
$ws.Range("A51").Copy()
$ws.Range("A52").PasteSpecial($xlPasteFormats)
$ws.Range("B52").Value = 5
$ws.Range("B51").Copy()
$ws.Range("B52").PasteSpecial($xlPasteFormats)
$ws.Range("C52").Value = 0.08333333333333333
$ws.Range("C51").Copy()
$ws.Range("C52").PasteSpecial($xlPasteFormats)
$ws.Range("D52").Value = "Configuration"
$ws.Range("D51").Copy()
$ws.Range("D52").PasteSpecial($xlPasteFormats)
$ws.Range("E52").Value = "Configuration des catégories pour le service Développement"

# --- Fix up the old placeholder task: it becomes the write-up for the
#     support IT interview report. ---
$ws.Range("E45").Value = "Rédaction d'un rapport d'entretien "

# --- Expand table "Tableau1" + autofilter to cover the new rows ---
$lo = $ws.ListObjects.Item("Tableau1")
$lo.Resize($ws.Range("A1:F52"))

# --- Update selections to match the post-edit cursor position ---
$ws.Application.Goto($ws.Range("D53"))
$totaux.Application.Goto($totaux.Range("I7"))
$ws.Application.Goto($ws.Range("D53"))

# Re-select the Journal sheet view state (topLeftCell stays the same)
$ws.Activate()
